$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J: header "Feedbacks dos Paineis" ---
$ws.Range("J1").Value = "Feedbacks dos Paineis"
# Reuse the exact header formatting (bold, centered, bordered) from the
# existing header row by copying format only (keeps same style slot).
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- New row 10 data ---
$ws.Range("A10").Value = "cleiton.souza@mrv.com.br"
$ws.Range("B10").Value = "Planilha automatizada"
$ws.Range("C10").Value = "Automação"
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = "Muito importante `n"
$ws.Range("F10").Value = "Painel Produção Produtividade e MO - PLNESROBR005; PAP - Dossiê; Painel Operações - Planejamento e Controle - PLNESROBR010"
$ws.Range("I10").Value = "2025-05-19 19:28:59"
$ws.Range("J10").Value = "Painel Produção Produtividade e MO - PLNESROBR005: gsadgfasdgasg; PAP - Dossiê: asdfghjkl`n"
